$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.406.92"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.585.15"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.46"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.78"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").Value = "3.583.60"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.92"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "4.191.61"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000204"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.51"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("D16").Value = "3.585.19"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "66.474.38"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.16"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.21"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "3.732.65"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.34"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.49"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "3.583.59"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.90"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.41"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.72"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  -4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "174.80"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0849"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.18"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.879"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.05"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.84"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.82"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.26"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.11"
$ws.Range("E51").Value = "  -0.41%  "
